# Applies the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" (sheet1): update Version/Date/Publisher, add Jurisdiction, drop Contact ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# The "Contact" / "No display for ContactDetail" row (row 10) becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The duplicate "Contact" row (row 11) is removed entirely, shifting remaining rows up
$meta.Rows.Item(11).Delete()

# ---- Sheet "Elements" (sheet2): update root Extension row's Short/Definition text ----
$elem = $wb.Worksheets.Item("Elements")

$elem.Range("K2").Value = "Episode Use Scale"
$elem.Range("L2").Value = "Use scale for the episode of care"
